# Add a new parameter row "general.maxNumberCompThreads" to the functional
# parameter properties sheet, inserted as the new row 18 (pushing the
# existing rows 18-53 down to 19-54).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 18 - this shifts rows 18..53 down to 19..54.
$insertRow = $ws.Rows.Item(18)
$insertRow.Insert()

# Fill in the values for the newly inserted row 18.
$ws.Range("A18").Value = "general.maxNumberCompThreads"
$ws.Range("D18").Value = "structural_pipeline"
$ws.Range("E18").Value = "numeric"
$ws.Range("F18").Value = "scalar nonempty nonnegative"
$ws.Range("G18").Value = "standard"
$ws.Range("H18").Value = "Maximum number of computational threads used in pipeline. Value 0 lets MATLAB determine the most desirable number of computational threads (equal to the number of physical cores on the machine)."

# Match the style used by the other rows' F and G columns (text-formatted).
$ws.Range("F18").NumberFormat = "@"
$ws.Range("G18").NumberFormat = "@"

# Move the active selection, matching the recorded end-user cursor position.
$ws.Range("A21").Select()
